# Update "想去人数" (attendance) figures for the refreshed data pull.
$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 7989
$ws1.Range("F5").Value = 5832
$ws1.Range("F11").Value = 359

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 89

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 7989
$ws4.Range("F5").Value = 5832
$ws4.Range("F11").Value = 89
$ws4.Range("F14").Value = 359
